$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Menu"
$ws.Cells.Item(2, 2).Value = 15
$ws.Cells.Item(2, 3).Value = 5
$ws.Cells.Item(2, 4).Value = 20

$ws.Cells.Item(3, 1).Value = "Tree"
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 10

$ws.Cells.Item(4, 1).Value = "DataGrid"
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 4).Value = 9

$ws.Cells.Item(5, 1).Value = "Popover"
$ws.Cells.Item(5, 2).Value = 7
$ws.Cells.Item(5, 3).Value = 2
$ws.Cells.Item(5, 4).Value = 9

$ws.Cells.Item(6, 1).Value = "Nav"
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 9

$ws.Cells.Item(7, 1).Value = "Dialog"
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 7

$ws.Cells.Item(8, 1).Value = "Combobox"
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 6

$ws.Cells.Item(9, 1).Value = "TagPicker"
$ws.Cells.Item(9, 2).Value = 4
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 5

$ws.Cells.Item(10, 1).Value = "Tooltip"
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 5

$ws.Cells.Item(11, 1).Value = "Table"
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 5

$ws.Cells.Item(12, 1).Value = "Toolbar"
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Value = 4

$ws.Cells.Item(13, 1).Value = "Virtualizer"
$ws.Cells.Item(13, 2).Value = 4
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 4

$ws.Cells.Item(14, 1).Value = "Dropdown"
$ws.Cells.Item(14, 2).Value = 2
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = 4

$ws.Cells.Item(15, 1).Value = "MessageBar"
$ws.Cells.Item(15, 2).Value = 3
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 3

$ws.Cells.Item(16, 1).Value = "TeachingPopover"
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 3

$ws.Cells.Item(17, 1).Value = "Skeleton"
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = 3

$ws.Cells.Item(18, 1).Value = "Calendar Compat"
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 3

$ws.Cells.Item(19, 1).Value = "Accordion"
$ws.Cells.Item(19, 2).Value = 2
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 2

$ws.Cells.Item(20, 1).Value = "DatePicker"
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 2

$ws.Cells.Item(21, 1).Value = "Drawer"
$ws.Cells.Item(21, 2).Value = 2
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 2

$ws.Cells.Item(22, 1).Value = "Slider"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = 2

$ws.Cells.Item(23, 1).Value = "Tabs"
$ws.Cells.Item(23, 2).Value = 2
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 2

$ws.Cells.Item(24, 1).Value = "List"
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 2

$ws.Cells.Item(25, 1).Value = "FluentProvider"
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = 2

$ws.Cells.Item(26, 1).Value = "DatePickerCompat"
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 2

$ws.Cells.Item(27, 1).Value = "Switch"
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 2

$ws.Cells.Item(28, 1).Value = "Portal"
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 2

$ws.Cells.Item(29, 1).Value = "Toast"
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 2

$ws.Cells.Item(30, 1).Value = "Input"
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 1

$ws.Cells.Item(31, 1).Value = "Label"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 1

$ws.Cells.Item(32, 1).Value = "Popup"
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 1

$ws.Cells.Item(33, 1).Value = "SwatchPicker"
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 1

$ws.Cells.Item(34, 1).Value = "FocusTrapZone"
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 1

$ws.Cells.Item(35, 1).Value = "AvatarGroup"
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 1

$ws.Cells.Item(36, 1).Value = "Avatar"
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 1

$ws.Cells.Item(37, 1).Value = "MenuItem"
$ws.Cells.Item(37, 2).Value = 1
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 1

$ws.Cells.Item(38, 1).Value = "PopoverSurface"
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 1

$ws.Cells.Item(39, 1).Value = "Checkbox"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 1

$ws.Cells.Item(40, 1).Value = "Button"
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = 1

$ws.Cells.Item(41, 1).Value = "Image"
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = 1

$ws.Cells.Item(42, 1).Value = "Badge"
$ws.Cells.Item(42, 2).Value = 1
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 1

$ws.Cells.Item(43, 1).Value = "TabList"
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 1
$ws.Cells.Item(43, 4).Value = 1

$ws.Cells.Item(44, 1).Value = "Tag"
$ws.Cells.Item(44, 2).Value = 1
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 1

$ws.Cells.Item(45, 1).Value = "Rating"
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 0

$ws.Cells.Item(46, 1).Value = "InfoLabel"
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 0

$ws.Cells.Item(47, 1).Value = "Segment"
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 0

$ws.Cells.Item(48, 1).Value = "ColorPicker"
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 0

$ws.Cells.Item(49, 1).Value = "SplitButton"
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 0

$ws.Cells.Item(50, 1).Value = "Pickers"
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 0

$ws.Cells.Item(51, 1).Value = "Keytip"
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 0

$ws.Cells.Item(52, 1).Value = "SpinButton"
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 0
